$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new log entry as row 30
$ws.Range("A30").Value = "Afmelding nieuwsbrief"
$ws.Range("B30").Value = "mailmind.test@zohomail.eu"
$ws.Range("C30").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$ws.Range("D30").Value = "Afmelding"
$ws.Range("F30").Value = "2025-06-19 17:56:20"
$ws.Range("G30").Value = "Nee"

# Extend conditional formatting ranges to include the new row
$ws.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D30"))
$ws.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G30"))

# Update the Dashboard summary sheet
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A2").Value = "Afmelding"
$dash.Range("B2").Value = 8
$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 8
